$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.420.97"
$ws.Range("E2").Value = "  -3.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.20"
$ws.Range("E3").Value = "  -3.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.12"
$ws.Range("E6").Value = "  -2.67%  "

$ws.Range("E7").Value = "  -2.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.29"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3409"
$ws.Range("E9").Value = "  -1.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.173"
$ws.Range("E10").Value = "  -2.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07671"
$ws.Range("E11").Value = "  -4.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.34"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.070"
$ws.Range("E14").Value = "  -3.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.942"
$ws.Range("E15").Value = "  -3.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.582.31"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("E17").Value = "  -4.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.41"
$ws.Range("E18").Value = "  -4.32%  "

$ws.Range("E19").Value = "  -3.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.273"
$ws.Range("E21").Value = "  -4.85%  "

$ws.Range("E22").Value = "  -3.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5339"
$ws.Range("E23").Value = "  -6.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.03"
$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.428.27"
$ws.Range("E25").Value = "  -3.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.356"
$ws.Range("E26").Value = "  -2.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.881"
$ws.Range("E27").Value = "  -2.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.11"
$ws.Range("E28").Value = "  -3.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "146.93"
$ws.Range("E29").Value = "  -1.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.988"
$ws.Range("E30").Value = "  -3.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.88"
$ws.Range("E31").Value = "  -3.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.752.97"
$ws.Range("E32").Value = "  -2.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.023"
$ws.Range("E33").Value = "  +4.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.240"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.018"
$ws.Range("E35").Value = "  -4.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.10"
$ws.Range("E36").Value = "  -8.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08477"
$ws.Range("E37").Value = "  -2.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02549"
$ws.Range("E38").Value = "  -3.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2323"
$ws.Range("E39").Value = "  -3.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.554"
$ws.Range("E40").Value = "  -3.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06489"
$ws.Range("E41").Value = "  -3.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.301"
$ws.Range("E42").Value = "  +1.42%  "

$ws.Range("E43").Value = "  -7.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6366"
$ws.Range("E44").Value = "  -5.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.23"
$ws.Range("E45").Value = "  -6.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9996"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6002"
$ws.Range("E47").Value = "  -4.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.761"
$ws.Range("E48").Value = "  -3.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.115"
$ws.Range("E49").Value = "  -4.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.262"
$ws.Range("E50").Value = "  +3.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.74"
$ws.Range("E51").Value = "  -1.03%  "
